# Update "sweep data" worksheet rows 2-22 (the 7 "Step Diameters" groups x 3
# "Thrust" rows each) with re-computed path-planning results.
#
# Columns updated for every data row (2-22):
#   G  Total Rocket Mass        -> constant across all rows now
#   H  Rocket Height            -> constant per Step-Diameter group (7 groups of 3 rows)
#   I  Stage 1 Mass             -> constant across all rows now
#   J  Stage 2 Mass             -> constant across all rows now
#   K  Stage 1 Structural Factor-> constant across all rows now
#   L  Stage 2 Structural Factor-> constant across all rows now
#   O  Engine count (step 1)    -> cycles 10, 9, 8 with the Thrust sub-row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that are now identical for every data row (2-22)
$totalRocketMass   = 1458377.232148911
$stage1Mass        = 1344745.458720702
$stage2Mass        = 113631.7734282086
$stage1StructFactor = 0.0461
$stage2StructFactor = 0.101

# Rocket Height (column H) is constant within each block of 3 rows that
# share the same "Step Diameters" value (rows 2-4, 5-7, 8-10, ... 20-22)
$rocketHeightByGroup = @(
    64.11190241563072,
    61.28830538310085,
    58.77510906724572,
    56.53291518611289,
    54.5283836203006,
    52.73314798945927,
    51.12295111925128
)

# Engine count (step 1), column O, cycles through these 3 values for every
# group of 3 rows (one value per distinct Thrust within the group)
$engineCountCycle = @(10, 9, 8)

$row = 2
for ($group = 0; $group -lt 7; $group++) {
    $rocketHeight = $rocketHeightByGroup[$group]
    for ($sub = 0; $sub -lt 3; $sub++) {
        $ws.Range("G$row").Value = $totalRocketMass
        $ws.Range("H$row").Value = $rocketHeight
        $ws.Range("I$row").Value = $stage1Mass
        $ws.Range("J$row").Value = $stage2Mass
        $ws.Range("K$row").Value = $stage1StructFactor
        $ws.Range("L$row").Value = $stage2StructFactor
        $ws.Range("O$row").Value = $engineCountCycle[$sub]

        $row++
    }
}

Write-Host "Updated rows 2-22 of sweep data (G, H, I, J, K, L, O columns)."
